# Reorder "Recorded By" (column G) values so that "System" appears first
# in the comma-separated list, e.g. "dnasr281@gmail.com, System" becomes
# "System, dnasr281@gmail.com". Entries that already start with "System"
# (case-sensitive) or that do not contain a "System" token after the first
# position are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    # find index of the literal "System" token among parts after the first
    $sysIndex = -1
    for ($i = 1; $i -lt $parts.Count; $i++) {
        if ($parts[$i] -ceq "System") {
            $sysIndex = $i
            break
        }
    }

    if ($sysIndex -lt 0) { continue }
    if ($parts[0] -ceq "System") { continue }

    $newParts = @("System")
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -eq $sysIndex) { continue }
        $newParts += $parts[$i]
    }

    $cell.Value = [string]::Join(", ", $newParts)
}
